$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.719623804092407
$ws.Range("B1").Value = 2.106775522232056
$ws.Range("C1").Value = 2.059844017028809
$ws.Range("D1").Value = 1.906466603279114
$ws.Range("E1").Value = 1.547598600387573
